$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (Q5) updated values
$ws.Range("B7").Value = 0.2555986832542369
$ws.Range("C7").Value = 0.4969596607080333
$ws.Range("D7").Value = 0.4509468104305381
$ws.Range("E7").Value = 0.6715257332601172
$ws.Range("F7").Value = 0.658648722000501
$ws.Range("G7").Value = 9

# Row 8 (Q6) updated values
$ws.Range("B8").Value = 0.1981259026106764
$ws.Range("C8").Value = 0.4731767774008284
$ws.Range("D8").Value = 0.4677985065100134
$ws.Range("E8").Value = 0.6839579713038026
$ws.Range("F8").Value = 0.6943433677783694
$ws.Range("G8").Value = 9
